$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '23.453.55'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.29%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.647.75'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.35%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.000'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '298.44'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.85%  '
$ws.Range('E7').Value = '  -1.29%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3555'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.65%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '49.94'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.65%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08096'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.76%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.219'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.60%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9999'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.03%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.07'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.88%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.394'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.43%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.354'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.88%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001197'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.20%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.646.65'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.40%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '97.30'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.49%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06947'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.53%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.758'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.32%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.32'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.29%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9996'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.43'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.83%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '23.467.19'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.18%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.494'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.63%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.900'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -5.76%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.89'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.10%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '152.93'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.15%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.208'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.69%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '132.78'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.65%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.835.63'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.02%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.922'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.70%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.123'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.82%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.45'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.24%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9965'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -8.34%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02715'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.44%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.08733'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.22%  '
$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.935'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.64%  '
$ws.Range('B39').Value = 'Algorand'
$ws.Range('C39').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2428'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.67%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '13.05'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.75%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.06771'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.99%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6879'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.88%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.304'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.85%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '15.51'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.70%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.000'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.12%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6358'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.20%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.256'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.55%  '
$ws.Range('E48').Value = '  -1.79%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.07728'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.31%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '127.21'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.65%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.150'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.09%  '
